$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9223046214701632
$ws.Range("J2").Value = 0.9223046214701632
$ws.Range("M2").Value = 135.0916853333333
$ws.Range("N2").Value = 405.2750559999999
$ws.Range("O2").Value = 0.7123704212620513
$ws.Range("P2").Value = 0.7123704212620514
$ws.Range("Q2").Value = 225.4033139040586
$ws.Range("R2").Value = 2028.629825136528
$ws.Range("S2").Value = 0.6570225317286369
$ws.Range("T2").Value = 0.657022531728637
$ws.Range("I3").Value = 0.9223046214701632
$ws.Range("J3").Value = 0.9223046214701632
$ws.Range("O3").Value = 0.2125756143240238
$ws.Range("P3").Value = 0.2125756143240238
$ws.Range("S3").Value = 0.1960594715029061
$ws.Range("T3").Value = 0.1960594715029062
$ws.Range("I4").Value = 0.9223046214701632
$ws.Range("J4").Value = 0.9223046214701632
$ws.Range("M4").Value = 14.23299766666667
$ws.Range("N4").Value = 42.698993
$ws.Range("O4").Value = 0.07505396441392481
$ws.Range("P4").Value = 0.07505396441392483
$ws.Range("Q4").Value = 23.74805549978434
$ws.Range("R4").Value = 213.732499498059
$ws.Range("S4").Value = 0.06922261823862003
$ws.Range("T4").Value = 0.06922261823862004
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.140557
$ws.Range("H5").Value = 0.421671
$ws.Range("I5").Value = 0.07769537852983674
$ws.Range("J5").Value = 0.07769537852983674
$ws.Range("M5").Value = 135.0916853333333
$ws.Range("N5").Value = 405.2750559999999
$ws.Range("O5").Value = 0.7123704212620513
$ws.Range("P5").Value = 0.7123704212620514
$ws.Range("Q5").Value = 18.98808201539733
$ws.Range("R5").Value = 170.892738138576
$ws.Range("S5").Value = 0.05534788953341433
$ws.Range("T5").Value = 0.05534788953341434
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.140557
$ws.Range("H6").Value = 0.421671
$ws.Range("I6").Value = 0.07769537852983674
$ws.Range("J6").Value = 0.07769537852983674
$ws.Range("O6").Value = 0.2125756143240238
$ws.Range("P6").Value = 0.2125756143240238
$ws.Range("Q6").Value = 5.666157772394667
$ws.Range("R6").Value = 50.995419951552
$ws.Range("S6").Value = 0.01651614282111761
$ws.Range("T6").Value = 0.01651614282111762
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.140557
$ws.Range("H7").Value = 0.421671
$ws.Range("I7").Value = 0.07769537852983674
$ws.Range("J7").Value = 0.07769537852983674
$ws.Range("M7").Value = 14.23299766666667
$ws.Range("N7").Value = 42.698993
$ws.Range("O7").Value = 0.07505396441392481
$ws.Range("P7").Value = 0.07505396441392483
$ws.Range("Q7").Value = 2.000547453033667
$ws.Range("R7").Value = 18.004927077303
$ws.Range("S7").Value = 0.005831346175304785
$ws.Range("T7").Value = 0.005831346175304786
